$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.04742086062181
$ws.Cells.Item(2, 4).Value = 1.047682729875347
$ws.Cells.Item(2, 5).Value = 1.054428496499307
$ws.Cells.Item(2, 6).Value = 1.064124011899746
$ws.Cells.Item(2, 9).Value = 1.046784447897723
$ws.Cells.Item(2, 10).Value = 1.052470020235141
$ws.Cells.Item(2, 11).Value = 1.050444683892116
$ws.Cells.Item(2, 12).Value = 1.057171733440442
$ws.Cells.Item(2, 13).Value = 1.066840797224516
$ws.Cells.Item(2, 14).Value = 1.021256813391929

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.048575128183609
$ws.Cells.Item(3, 4).Value = 1.04835470381929
$ws.Cells.Item(3, 5).Value = 1.05546374015635
$ws.Cells.Item(3, 6).Value = 1.065262841562896
$ws.Cells.Item(3, 9).Value = 1.047107889716463
$ws.Cells.Item(3, 10).Value = 1.053271747076397
$ws.Cells.Item(3, 11).Value = 1.050928914006978
$ws.Cells.Item(3, 12).Value = 1.058019650951801
$ws.Cells.Item(3, 13).Value = 1.067793959835154
$ws.Cells.Item(3, 14).Value = 1.02153024344587

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.049321926399497
$ws.Cells.Item(4, 4).Value = 1.048789227558983
$ws.Cells.Item(4, 5).Value = 1.056133861673451
$ws.Cells.Item(4, 6).Value = 1.066000116533871
$ws.Cells.Item(4, 9).Value = 1.047315733170069
$ws.Cells.Item(4, 10).Value = 1.053789880939701
$ws.Cells.Item(4, 11).Value = 1.051241274116438
$ws.Cells.Item(4, 12).Value = 1.058567953980018
$ws.Cells.Item(4, 13).Value = 1.068410508649997
$ws.Cells.Item(4, 14).Value = 1.02170678621664

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.049635859937461
$ws.Cells.Item(5, 4).Value = 1.048971831544739
$ws.Cells.Item(5, 5).Value = 1.05641564091511
$ws.Cells.Item(5, 6).Value = 1.066310157472909
$ws.Cells.Item(5, 9).Value = 1.047402764575952
$ws.Cells.Item(5, 10).Value = 1.05400755237561
$ws.Cells.Item(5, 11).Value = 1.051372357894704
$ws.Cells.Item(5, 12).Value = 1.058798375683013
$ws.Cells.Item(5, 13).Value = 1.06866965567364
$ws.Cells.Item(5, 14).Value = 1.021780912738851

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.049688569594716
$ws.Cells.Item(6, 4).Value = 1.04900248744228
$ws.Cells.Item(6, 5).Value = 1.05646295642859
$ws.Cells.Item(6, 6).Value = 1.066362220035322
$ws.Cells.Item(6, 9).Value = 1.047417357242262
$ws.Cells.Item(6, 10).Value = 1.05404409146522
$ws.Cells.Item(6, 11).Value = 1.051394353797608
$ws.Cells.Item(6, 12).Value = 1.058837059534752
$ws.Cells.Item(6, 13).Value = 1.068713164674435
$ws.Cells.Item(6, 14).Value = 1.021793353511289

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.049326121273655
$ws.Cells.Item(7, 4).Value = 1.048791667797597
$ws.Cells.Item(7, 5).Value = 1.056137626584609
$ws.Cells.Item(7, 6).Value = 1.066004258959656
$ws.Cells.Item(7, 9).Value = 1.047316897447499
$ws.Cells.Item(7, 10).Value = 1.053792790072919
$ws.Cells.Item(7, 11).Value = 1.051243026578021
$ws.Cells.Item(7, 12).Value = 1.058571033219141
$ws.Cells.Item(7, 13).Value = 1.068413971581504
$ws.Cells.Item(7, 14).Value = 1.021707777060518

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.047810970082252
$ws.Cells.Item(8, 4).Value = 1.047909885858933
$ws.Cells.Item(8, 5).Value = 1.054778310484578
$ws.Cells.Item(8, 6).Value = 1.064508807170462
$ws.Cells.Item(8, 9).Value = 1.046894055755357
$ws.Cells.Item(8, 10).Value = 1.052741099751065
$ws.Cells.Item(8, 11).Value = 1.05060853216073
$ws.Cells.Item(8, 12).Value = 1.057458365154212
$ws.Cells.Item(8, 13).Value = 1.067162966976887
$ws.Cells.Item(8, 14).Value = 1.02134930005775

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.045140339139366
$ws.Cells.Item(9, 4).Value = 1.04635389842871
$ws.Cells.Item(9, 5).Value = 1.052384915845204
$ws.Cells.Item(9, 6).Value = 1.061876486195079
$ws.Cells.Item(9, 9).Value = 1.046137888779833
$ws.Cells.Item(9, 10).Value = 1.050882987150078
$ws.Cells.Item(9, 11).Value = 1.04948306631263
$ws.Cells.Item(9, 12).Value = 1.055494954269451
$ws.Cells.Item(9, 13).Value = 1.0649568995161
$ws.Cells.Item(9, 14).Value = 1.020714668410616

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.043359346374218
$ws.Cells.Item(10, 4).Value = 1.045315154211427
$ws.Cells.Item(10, 5).Value = 1.050790564780289
$ws.Cells.Item(10, 6).Value = 1.060123489574345
$ws.Cells.Item(10, 9).Value = 1.045626337447466
$ws.Cells.Item(10, 10).Value = 1.049640919015335
$ws.Cells.Item(10, 11).Value = 1.04872779768237
$ws.Cells.Item(10, 12).Value = 1.05418413760006
$ws.Cells.Item(10, 13).Value = 1.063485055461225
$ws.Cells.Item(10, 14).Value = 1.020289593034575

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.042588003868028
$ws.Cells.Item(11, 4).Value = 1.044865038363987
$ws.Cells.Item(11, 5).Value = 1.05010048139519
$ws.Cells.Item(11, 6).Value = 1.059364860526004
$ws.Cells.Item(11, 9).Value = 1.045403064840874
$ws.Cells.Item(11, 10).Value = 1.049102292924089
$ws.Cells.Item(11, 11).Value = 1.048399585894056
$ws.Cells.Item(11, 12).Value = 1.053616087334093
$ws.Cells.Item(11, 13).Value = 1.062847454271467
$ws.Cells.Item(11, 14).Value = 1.020105057917146

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.042301467154079
$ws.Cells.Item(12, 4).Value = 1.044697795921619
$ws.Cells.Item(12, 5).Value = 1.049844195218853
$ws.Cells.Item(12, 6).Value = 1.059083135644644
$ws.Cells.Item(12, 9).Value = 1.045319865778467
$ws.Cells.Item(12, 10).Value = 1.048902101882119
$ws.Cells.Item(12, 11).Value = 1.048277497051893
$ws.Cells.Item(12, 12).Value = 1.053405018646399
$ws.Cells.Item(12, 13).Value = 1.062610577499583
$ws.Cells.Item(12, 14).Value = 1.020036441885975

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.042362931443302
$ws.Cells.Item(13, 4).Value = 1.044733672217735
$ws.Cells.Item(13, 5).Value = 1.049899167606439
$ws.Cells.Item(13, 6).Value = 1.059143563718588
$ws.Cells.Item(13, 9).Value = 1.045337724284202
$ws.Cells.Item(13, 10).Value = 1.048945049043124
$ws.Cells.Item(13, 11).Value = 1.048303693500405
$ws.Cells.Item(13, 12).Value = 1.053450296754159
$ws.Cells.Item(13, 13).Value = 1.062661390343299
$ws.Cells.Item(13, 14).Value = 1.020051163494388

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.042564319188635
$ws.Cells.Item(14, 4).Value = 1.044851215055228
$ws.Cells.Item(14, 5).Value = 1.050079295857106
$ws.Cells.Item(14, 6).Value = 1.059341571764583
$ws.Cells.Item(14, 9).Value = 1.045396193000855
$ws.Cells.Item(14, 10).Value = 1.049085747549703
$ws.Cells.Item(14, 11).Value = 1.048389497588844
$ws.Cells.Item(14, 12).Value = 1.053598641752998
$ws.Cells.Item(14, 13).Value = 1.062827874851468
$ws.Cells.Item(14, 14).Value = 1.020099387553528

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.042688397363681
$ws.Cells.Item(15, 4).Value = 1.044923630552433
$ws.Cells.Item(15, 5).Value = 1.050190284297395
$ws.Cells.Item(15, 6).Value = 1.05946357947448
$ws.Cells.Item(15, 9).Value = 1.045432182298526
$ws.Cells.Item(15, 10).Value = 1.049172420442599
$ws.Cells.Item(15, 11).Value = 1.04844234095575
$ws.Cells.Item(15, 12).Value = 1.053690032768368
$ws.Cells.Item(15, 13).Value = 1.062930445690401
$ws.Cells.Item(15, 14).Value = 1.020129090512483

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.043410534310463
$ws.Cells.Item(16, 4).Value = 1.045345019959963
$ws.Cells.Item(16, 5).Value = 1.050836369203837
$ws.Cells.Item(16, 6).Value = 1.060173846271019
$ws.Cells.Item(16, 9).Value = 1.045641118064898
$ws.Cells.Item(16, 10).Value = 1.049676648883154
$ws.Cells.Item(16, 11).Value = 1.048749555271665
$ws.Cells.Item(16, 12).Value = 1.054221827514691
$ws.Cells.Item(16, 13).Value = 1.063527364922993
$ws.Cells.Item(16, 14).Value = 1.020301830008279

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.043863467541881
$ws.Cells.Item(17, 4).Value = 1.04560925775389
$ws.Cells.Item(17, 5).Value = 1.051241715959939
$ws.Cells.Item(17, 6).Value = 1.060619492381436
$ws.Cells.Item(17, 9).Value = 1.045771704650667
$ws.Cells.Item(17, 10).Value = 1.049992722690143
$ws.Cells.Item(17, 11).Value = 1.048941948150378
$ws.Cells.Item(17, 12).Value = 1.054555285184539
$ws.Cells.Item(17, 13).Value = 1.0639017202679
$ws.Cells.Item(17, 14).Value = 1.020410057735449

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.044127640416417
$ws.Cells.Item(18, 4).Value = 1.045763351042371
$ws.Cells.Item(18, 5).Value = 1.051478175039938
$ws.Cells.Item(18, 6).Value = 1.060879471771691
$ws.Cells.Item(18, 9).Value = 1.045847703066942
$ws.Cells.Item(18, 10).Value = 1.050177005827336
$ws.Cells.Item(18, 11).Value = 1.049054054228606
$ws.Cells.Item(18, 12).Value = 1.054749741259766
$ws.Cells.Item(18, 13).Value = 1.064120048154124
$ws.Cells.Item(18, 14).Value = 1.020473139349047

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.044217713973319
$ws.Cells.Item(19, 4).Value = 1.045815887400216
$ws.Cells.Item(19, 5).Value = 1.051558806106805
$ws.Cells.Item(19, 6).Value = 1.060968125100453
$ws.Cells.Item(19, 9).Value = 1.045873587641954
$ws.Cells.Item(19, 10).Value = 1.050239828515716
$ws.Cells.Item(19, 11).Value = 1.04909226025685
$ws.Cells.Item(19, 12).Value = 1.054816038290779
$ws.Cells.Item(19, 13).Value = 1.064194487722455
$ws.Cells.Item(19, 14).Value = 1.020494640787354

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.043814873698511
$ws.Cells.Item(20, 4).Value = 1.045580910859815
$ws.Cells.Item(20, 5).Value = 1.051198223283355
$ws.Cells.Item(20, 6).Value = 1.060571674479444
$ws.Cells.Item(20, 9).Value = 1.045757711585418
$ws.Cells.Item(20, 10).Value = 1.049958818969827
$ws.Cells.Item(20, 11).Value = 1.048921317946409
$ws.Cells.Item(20, 12).Value = 1.054519512903205
$ws.Cells.Item(20, 13).Value = 1.063861558306628
$ws.Cells.Item(20, 14).Value = 1.020398450661796

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.042505016241069
$ws.Cells.Item(21, 4).Value = 1.044816602977347
$ws.Cells.Item(21, 5).Value = 1.050026251436536
$ws.Cells.Item(21, 6).Value = 1.059283261578258
$ws.Cells.Item(21, 9).Value = 1.045378982756061
$ws.Cells.Item(21, 10).Value = 1.049044318701882
$ws.Cells.Item(21, 11).Value = 1.048364235290241
$ws.Cells.Item(21, 12).Value = 1.053554959773283
$ws.Cells.Item(21, 13).Value = 1.062778850514942
$ws.Cells.Item(21, 14).Value = 1.020085188744314

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.041681307023302
$ws.Cells.Item(22, 4).Value = 1.044335766906729
$ws.Cells.Item(22, 5).Value = 1.049289626078082
$ws.Cells.Item(22, 6).Value = 1.058473553504301
$ws.Cells.Item(22, 9).Value = 1.045139323730915
$ws.Cells.Item(22, 10).Value = 1.048468633725878
$ws.Cells.Item(22, 11).Value = 1.048012955001823
$ws.Cells.Item(22, 12).Value = 1.052948104322394
$ws.Cells.Item(22, 13).Value = 1.062097858239547
$ws.Cells.Item(22, 14).Value = 1.019887814949242

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.042117985600371
$ws.Cells.Item(23, 4).Value = 1.044590694051275
$ws.Cells.Item(23, 5).Value = 1.049680102623566
$ws.Cells.Item(23, 6).Value = 1.058902760464788
$ws.Cells.Item(23, 9).Value = 1.045266517264441
$ws.Cells.Item(23, 10).Value = 1.048773882060394
$ws.Cells.Item(23, 11).Value = 1.048199271966157
$ws.Cells.Item(23, 12).Value = 1.053269848282374
$ws.Cells.Item(23, 13).Value = 1.062458889175416
$ws.Cells.Item(23, 14).Value = 1.019992485789249

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.043836831211808
$ws.Cells.Item(24, 4).Value = 1.045593719700386
$ws.Cells.Item(24, 5).Value = 1.051217875667538
$ws.Cells.Item(24, 6).Value = 1.060593281202342
$ws.Cells.Item(24, 9).Value = 1.04576403497608
$ws.Cells.Item(24, 10).Value = 1.04997413884044
$ws.Cells.Item(24, 11).Value = 1.048930640197992
$ws.Cells.Item(24, 12).Value = 1.054535676993463
$ws.Cells.Item(24, 13).Value = 1.06387970585209
$ws.Cells.Item(24, 14).Value = 1.02040369554016

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.045830856672921
$ws.Cells.Item(25, 4).Value = 1.046756411759193
$ws.Cells.Item(25, 5).Value = 1.053003444159809
$ws.Cells.Item(25, 6).Value = 1.062556669523207
$ws.Cells.Item(25, 9).Value = 1.046334687418887
$ws.Cells.Item(25, 10).Value = 1.051363937744269
$ws.Cells.Item(25, 11).Value = 1.049774901628926
$ws.Cells.Item(25, 12).Value = 1.056002871603178
$ws.Cells.Item(25, 13).Value = 1.065527418116848
$ws.Cells.Item(25, 14).Value = 1.020879085891637
